# Delete the "Time" / "Student Scheduler" row (row 2) from the schedule grid.
# This shifts every row below it up by one, which matches the target diff:
# the A-column index keeps its original value, and the B..G contents of the
# old rows 3-11 become the new rows 2-10.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Delete()
